$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 125.666664
$ws.Cells.Item(11, 9).Value = 125.666664
$ws.Cells.Item(11, 11).Value = 125.666664
$ws.Cells.Item(11, 13).Value = 14.333336

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 850.7143
$ws.Cells.Item(17, 10).Value = 850.7143
$ws.Cells.Item(17, 12).Value = 2552.1429
$ws.Cells.Item(17, 14).Value = -2888.1429

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 3875.7778
$ws.Cells.Item(29, 9).Value = 1435.5
$ws.Cells.Item(29, 10).Value = 5828
$ws.Cells.Item(29, 11).Value = 4306.5
$ws.Cells.Item(29, 12).Value = 17484
$ws.Cells.Item(29, 13).Value = -4025.5
$ws.Cells.Item(29, 14).Value = -18046

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 3001.5
$ws.Cells.Item(74, 9).Value = 3001.5
$ws.Cells.Item(74, 11).Value = 3001.5
$ws.Cells.Item(74, 13).Value = -2065.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(77, 8).Value = 3001.5
$ws.Cells.Item(77, 9).Value = 3001.5
$ws.Cells.Item(77, 11).Value = 15007.5
$ws.Cells.Item(77, 13).Value = -10327.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(95, 8).Value = 14959
$ws.Cells.Item(95, 10).Value = 14959
$ws.Cells.Item(95, 12).Value = 14959
$ws.Cells.Item(95, 14).Value = -20451

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(140, 8).Value = 34890
$ws.Cells.Item(140, 10).Value = 34890
$ws.Cells.Item(140, 12).Value = 34890
$ws.Cells.Item(140, 14).Value = -45250

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3587
$ws.Cells.Item(32, 9).Value = 3710.318
$ws.Cells.Item(32, 10).Value = 3285.5557
$ws.Cells.Item(32, 11).Value = 3710.318
$ws.Cells.Item(32, 12).Value = 3285.5557
$ws.Cells.Item(32, 13).Value = -3423.318
$ws.Cells.Item(32, 14).Value = -3859.5557

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 5869.154
$ws.Cells.Item(61, 9).Value = 4724.875
$ws.Cells.Item(61, 11).Value = 4724.875
$ws.Cells.Item(61, 13).Value = -4512.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1731.4546
$ws.Cells.Item(74, 9).Value = 1731.4546
$ws.Cells.Item(74, 11).Value = 1731.4546
$ws.Cells.Item(74, 13).Value = -857.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(76, 8).Value = 18562.375
$ws.Cells.Item(76, 10).Value = 18562.375
$ws.Cells.Item(76, 12).Value = 18562.375
$ws.Cells.Item(76, 14).Value = -19238.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 1731.4546
$ws.Cells.Item(77, 9).Value = 1731.4546
$ws.Cells.Item(77, 11).Value = 8657.273000000001
$ws.Cells.Item(77, 13).Value = -4289.273000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(79, 8).Value = 18562.375
$ws.Cells.Item(79, 10).Value = 18562.375
$ws.Cells.Item(79, 12).Value = 18562.375
$ws.Cells.Item(79, 14).Value = -20902.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 5869.154
$ws.Cells.Item(136, 9).Value = 4724.875
$ws.Cells.Item(136, 11).Value = 14174.625
$ws.Cells.Item(136, 13).Value = -11624.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(12, 8).Value = 673
$ws.Cells.Item(12, 9).Value = 673
$ws.Cells.Item(12, 11).Value = 673
$ws.Cells.Item(12, 13).Value = -505

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1915.4286
$ws.Cells.Item(20, 9).Value = 1324.75
$ws.Cells.Item(20, 10).Value = 2703
$ws.Cells.Item(20, 11).Value = 1324.75
$ws.Cells.Item(20, 12).Value = 2703
$ws.Cells.Item(20, 13).Value = -1077.75
$ws.Cells.Item(20, 14).Value = -3197

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 745.6
$ws.Cells.Item(22, 10).Value = 1500
$ws.Cells.Item(22, 12).Value = 1500
$ws.Cells.Item(22, 14).Value = -1846

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3462.8462
$ws.Cells.Item(134, 9).Value = 3462.8462
$ws.Cells.Item(134, 11).Value = 10388.5386
$ws.Cells.Item(134, 13).Value = -7853.5386

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(19, 8).Value = 210.25
$ws.Cells.Item(19, 9).Value = 168.63637
$ws.Cells.Item(19, 10).Value = 301.8
$ws.Cells.Item(19, 11).Value = 168.63637
$ws.Cells.Item(19, 12).Value = 301.8
$ws.Cells.Item(19, 13).Value = 1.363630000000001
$ws.Cells.Item(19, 14).Value = -641.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(24, 8).Value = 210.25
$ws.Cells.Item(24, 9).Value = 168.63637
$ws.Cells.Item(24, 10).Value = 301.8
$ws.Cells.Item(24, 11).Value = 168.63637
$ws.Cells.Item(24, 12).Value = 301.8
$ws.Cells.Item(24, 13).Value = 1.363630000000001
$ws.Cells.Item(24, 14).Value = -641.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 7825.375
$ws.Cells.Item(86, 9).Value = 7433.8335
$ws.Cells.Item(86, 11).Value = 7433.8335
$ws.Cells.Item(86, 13).Value = -6310.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(89, 8).Value = 7825.375
$ws.Cells.Item(89, 9).Value = 7433.8335
$ws.Cells.Item(89, 11).Value = 37169.1675
$ws.Cells.Item(89, 13).Value = -31553.1675

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(96, 8).Value = 7718.8335
$ws.Cells.Item(96, 10).Value = 7718.8335
$ws.Cells.Item(96, 12).Value = 7718.8335
$ws.Cells.Item(96, 14).Value = -13210.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 1493
$ws.Cells.Item(122, 9).Value = 1493
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 4479
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -2029
$ws.Cells.Item(122, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(130, 8).Value = 24995
$ws.Cells.Item(130, 10).Value = 24995
$ws.Cells.Item(130, 12).Value = 24995
$ws.Cells.Item(130, 14).Value = -35035

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 340265.3
$ws.Cells.Item(4, 9).Value = 378033.78
$ws.Cells.Item(4, 10).Value = 349
$ws.Cells.Item(4, 11).Value = 1134101.34
$ws.Cells.Item(4, 12).Value = 1047
$ws.Cells.Item(4, 13).Value = -1133989.34
$ws.Cells.Item(4, 14).Value = -1271

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 978.61536
$ws.Cells.Item(34, 10).Value = 1487.5
$ws.Cells.Item(34, 12).Value = 4462.5
$ws.Cells.Item(34, 14).Value = -4630.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(125, 8).Value = 8042.857
$ws.Cells.Item(125, 9).Value = 5433.3335
$ws.Cells.Item(125, 11).Value = 16300.0005
$ws.Cells.Item(125, 13).Value = -11380.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 1711.5238
$ws.Cells.Item(129, 9).Value = 986.75
$ws.Cells.Item(129, 10).Value = 2157.5386
$ws.Cells.Item(129, 11).Value = 2960.25
$ws.Cells.Item(129, 12).Value = 6472.6158
$ws.Cells.Item(129, 13).Value = 2039.75
$ws.Cells.Item(129, 14).Value = -16472.6158

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(101, 8).Value = 28664.666
$ws.Cells.Item(101, 10).Value = 28664.666
$ws.Cells.Item(101, 12).Value = 28664.666
$ws.Cells.Item(101, 14).Value = -35154.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 5749.8335
$ws.Cells.Item(132, 9).Value = 5749.8335
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 17249.5005
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -14719.5005
$ws.Cells.Item(132, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 1599
$ws.Cells.Item(40, 10).Value = 998
$ws.Cells.Item(40, 12).Value = 998
$ws.Cells.Item(40, 14).Value = -1270

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 4863.1177
$ws.Cells.Item(46, 9).Value = 4280.4443
$ws.Cells.Item(46, 10).Value = 5518.625
$ws.Cells.Item(46, 11).Value = 4280.4443
$ws.Cells.Item(46, 12).Value = 5518.625
$ws.Cells.Item(46, 13).Value = -4092.4443
$ws.Cells.Item(46, 14).Value = -5894.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(50, 8).Value = 60000
$ws.Cells.Item(50, 9).Value = 60000
$ws.Cells.Item(50, 11).Value = 60000
$ws.Cells.Item(50, 13).Value = -59363

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 6494
$ws.Cells.Item(82, 9).Value = 6494
$ws.Cells.Item(82, 11).Value = 6494
$ws.Cells.Item(82, 13).Value = -6133

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 8).Value = 6494
$ws.Cells.Item(85, 9).Value = 6494
$ws.Cells.Item(85, 11).Value = 6494
$ws.Cells.Item(85, 13).Value = -5246

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 6636
$ws.Cells.Item(122, 9).Value = 5999.6665
$ws.Cells.Item(122, 11).Value = 17998.9995
$ws.Cells.Item(122, 13).Value = -15548.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 8233
$ws.Cells.Item(132, 9).Value = 4754.8887
$ws.Cells.Item(132, 10).Value = 16058.75
$ws.Cells.Item(132, 11).Value = 14264.6661
$ws.Cells.Item(132, 12).Value = 48176.25
$ws.Cells.Item(132, 13).Value = -11734.6661
$ws.Cells.Item(132, 14).Value = -53236.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 11776
$ws.Cells.Item(62, 9).Value = 10402.667
$ws.Cells.Item(62, 10).Value = 12600
$ws.Cells.Item(62, 11).Value = 10402.667
$ws.Cells.Item(62, 12).Value = 12600
$ws.Cells.Item(62, 13).Value = -9778.666999999999
$ws.Cells.Item(62, 14).Value = -13848

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(65, 8).Value = 11776
$ws.Cells.Item(65, 9).Value = 10402.667
$ws.Cells.Item(65, 10).Value = 12600
$ws.Cells.Item(65, 11).Value = 52013.335
$ws.Cells.Item(65, 12).Value = 63000
$ws.Cells.Item(65, 13).Value = -48893.335
$ws.Cells.Item(65, 14).Value = -69240

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1900
$ws.Cells.Item(122, 9).Value = 1795.238
$ws.Cells.Item(122, 11).Value = 5385.714
$ws.Cells.Item(122, 13).Value = -2935.714
